$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3511
$ws1.Range("G2").Value = 75
$ws1.Range("F3").Value = 746
$ws1.Range("G3").Value = "不可售"
$ws1.Range("F5").Value = 7016
$ws1.Range("F6").Value = 2957
$ws1.Range("F7").Value = 54
$ws1.Range("F13").Value = 11
$ws1.Range("F16").Value = 13

# --- Sheet: 演出 ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 29

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3511
$ws4.Range("G2").Value = 75
$ws4.Range("F3").Value = 29
$ws4.Range("F4").Value = 746
$ws4.Range("G4").Value = "不可售"
$ws4.Range("F6").Value = 7016
$ws4.Range("F7").Value = 2957
$ws4.Range("F8").Value = 54
$ws4.Range("F14").Value = 11
$ws4.Range("F17").Value = 13
